$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.430.51"
$ws.Range("E2").Value = "  -2.09%  "
$ws.Range("D3").Value = "1.864.52"
$ws.Range("E3").Value = "  -2.18%  "
$ws.Range("D4").Value = "'1.007"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'319.91"
$ws.Range("E5").Value = "  -1.53%  "
$ws.Range("D6").Value = "'1.005"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").Value = "'0.4418"
$ws.Range("E7").Value = "  -3.99%  "
$ws.Range("D8").Value = "'0.3719"
$ws.Range("E8").Value = "  -2.42%  "
$ws.Range("D9").Value = "'0.07534"
$ws.Range("E9").Value = "  -2.35%  "
$ws.Range("D10").Value = "'0.9373"
$ws.Range("E10").Value = "  -3.79%  "
$ws.Range("D11").Value = "'21.27"
$ws.Range("E11").Value = "  -3.09%  "
$ws.Range("D12").Value = "1.914.48"
$ws.Range("E12").Value = "  -0.17%  "
$ws.Range("D13").Value = "'6.716"
$ws.Range("E13").Value = "  -3.22%  "
$ws.Range("D14").Value = "'5.453"
$ws.Range("E14").Value = "  -3.62%  "
$ws.Range("D15").Value = "'0.06869"
$ws.Range("E15").Value = "  -2.74%  "
$ws.Range("D16").Value = "'1.007"
$ws.Range("E16").Value = "  +0.16%  "
$ws.Range("D17").Value = "'82.12"
$ws.Range("E17").Value = "  -1.91%  "
$ws.Range("D18").Value = "'0.000009125"
$ws.Range("E18").Value = "  -3.84%  "
$ws.Range("D19").Value = "'1.005"
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").Value = "'16.01"
$ws.Range("E20").Value = "  -3.67%  "
$ws.Range("D21").Value = "28.422.49"
$ws.Range("E21").Value = "  -1.98%  "
$ws.Range("D22").Value = "'5.127"
$ws.Range("E22").Value = "  -3.40%  "
$ws.Range("D23").Value = "'10.73"
$ws.Range("E23").Value = "  -1.39%  "
$ws.Range("D24").Value = "2.118.51"
$ws.Range("E24").Value = "  -1.47%  "
$ws.Range("D25").Value = "'2.042"
$ws.Range("E25").Value = "  -2.64%  "
$ws.Range("D26").Value = "'154.94"
$ws.Range("E26").Value = "  -1.94%  "
$ws.Range("D27").Value = "'18.42"
$ws.Range("E27").Value = "  -3.52%  "
$ws.Range("D28").Value = "'5.374"
$ws.Range("E28").Value = "  -4.26%  "
$ws.Range("D29").Value = "'114.74"
$ws.Range("E29").Value = "  -2.11%  "
$ws.Range("D30").Value = "'1.743"
$ws.Range("E30").Value = "  -5.40%  "
$ws.Range("D31").Value = "'0.09139"
$ws.Range("E31").Value = "  -1.29%  "
$ws.Range("D32").Value = "'0.8073"
$ws.Range("E32").Value = "  -5.97%  "
$ws.Range("D33").Value = "'4.881"
$ws.Range("E33").Value = "  -4.02%  "
$ws.Range("D34").Value = "'1.164"
$ws.Range("E34").Value = "  -6.18%  "
$ws.Range("D35").Value = "'2.939"
$ws.Range("E35").Value = "  -1.53%  "
$ws.Range("D36").Value = "'1.005"
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("D37").Value = "'1.131"
$ws.Range("E37").Value = "  -0.86%  "
$ws.Range("D38").Value = "'0.05474"
$ws.Range("E38").Value = "  -3.44%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'3.044"
$ws.Range("E39").Value = "  +10.27%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.01954"
$ws.Range("E40").Value = "  -4.00%  "
$ws.Range("D41").Value = "'7.148"
$ws.Range("E41").Value = "  -3.45%  "
$ws.Range("D42").Value = "'0.5280"
$ws.Range("E42").Value = "  -3.79%  "
$ws.Range("D43").Value = "'0.1683"
$ws.Range("E43").Value = "  -3.93%  "
$ws.Range("D44").Value = "'8.854"
$ws.Range("E44").Value = "  -4.96%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "'2.068"
$ws.Range("E45").Value = "  -0.93%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").Value = "'0.06796"
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").Value = "'0.4910"
$ws.Range("E47").Value = "  -4.98%  "
$ws.Range("B48").Value = "PEPE"
$ws.Range("C48").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D48").Value = "'0.000002536"
$ws.Range("E48").Value = "  -4.19%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'10.64"
$ws.Range("E49").Value = "  -5.10%  "
$ws.Range("D50").Value = "'107.52"
$ws.Range("E50").Value = "  -2.33%  "
$ws.Range("D51").Value = "'1.682"
$ws.Range("E51").Value = "  -5.04%  "
